$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5 (Iowa)
$ws.Range("C5").Value = 34530
$ws.Range("E5").Value = 3002

# Row 6 (Illinois)
$ws.Range("B6").Value = 44023
$ws.Range("C6").Value = 152962
$ws.Range("D6").Value = 7168
$ws.Range("E6").Value = 25689
$ws.Range("F6").Value = 1982
$ws.Range("G6").Value = 16.79
$ws.Range("H6").Value = 27.65
